$wb = $excel.ActiveWorkbook

# --- Add the new "Quotes" worksheet as the last (3rd) tab ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$quotes = $wb.Worksheets.Add($null, $lastSheet)
$quotes.Name = "Quotes"

# --- Fill in the quotes table. Cells are written in the same order the ---
# --- original author typed them so new shared strings line up exactly. ---
$quotes.Range("B1").Value = "Quote"
$quotes.Range("A1").Value = "Auteur"
$quotes.Range("D1").Value = "Bron"
$quotes.Range("A2").Value = "Epictetus"
$quotes.Range("B2").Value = "Het is onmogelijk om te leren wat je denkt reeds te weten"
$quotes.Range("C1").Value = "Quote english"
$quotes.Range("C2").Value = "It is impossible for a man to learn what he thinks he already knows"
$quotes.Range("D2").Value = "Scheurkalender NewScientist 2017"
$quotes.Range("A3").Value = "Marvin Minsky"
$quotes.Range("C3").Value = "We rarely recognize how wonderful it is that a person can traverse an entire lifetime without making a single really serious mistake — like putting a fork in one's eye or using a window instead of a door."
$quotes.Range("D3").Value = "Scheurkalender NewScientist 2017"

# --- Column widths for the quote columns (B and C) ---
$quotes.Columns.Item(2).ColumnWidth = 53
$quotes.Columns.Item(3).ColumnWidth = 53

# --- Page orientation for the new sheet ---
$quotes.PageSetup.Orientation = 1

# --- Update the remembered selection on the other two sheets ---
$planning = $wb.Worksheets.Item("Planning")
$planning.Range("A26").Select()

$mensen = $wb.Worksheets.Item("Mensen in de loop")
$mensen.Range("A10").Select()

# --- Make Quotes the active / tab-selected sheet with its own selection ---
$quotes.Activate()
$quotes.Range("C20").Select()
